$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (all these columns are stored as text
# in the source data, so each cell is pre-formatted as Text before the write to
# avoid Excel auto-coercing numeric-looking strings into Number cells, which
# would silently drop formatting like trailing zeros, leading zeros, "%" and "," ).
$updates = [ordered]@{
    'D2' = '310.90'
    'E2' = '8.21%'
    'G2' = '9'
    'D3' = '32.52'
    'E3' = '10.15%'
    'G3' = '9'
    'D4' = '5.343'
    'E4' = '4.26%'
    'G4' = '9'
    'D5' = '0.07641'
    'E5' = '14.05%'
    'G5' = '9'
    'D6' = '7.900'
    'E6' = '7.63%'
    'G6' = '9'
    'D7' = '3.765'
    'E7' = '10.50%'
    'G7' = '9'
    'D8' = '1.602'
    'E8' = '17.94%'
    'G8' = '9'
    'D9' = '0.9176'
    'E9' = '-0.12%'
    'G9' = '9'
    'D10' = '0.01759'
    'E10' = '2,623.11%'
    'G10' = '9'
    'D11' = '0.1723'
    'E11' = '8.67%'
    'G11' = '9'
    'D12' = '0.07668'
    'E12' = '13.71%'
    'G12' = '9'
    'D13' = '0.08277'
    'E13' = '7.70%'
    'G13' = '9'
    'E14' = '3.41%'
    'G14' = '9'
    'D15' = '0.09913'
    'E15' = '10.43%'
    'G15' = '9'
    'D16' = '0.001515'
    'E16' = '-4.18%'
    'G16' = '9'
    'D17' = '0.04574'
    'E17' = '1.71%'
    'G17' = '9'
    'D18' = '0.006166'
    'E18' = '-2.05%'
    'G18' = '9'
    'D19' = '3.476'
    'E19' = '0.61%'
    'G19' = '9'
    'D20' = '2.248'
    'E20' = '1.21%'
    'G20' = '9'
    'D21' = '0.3291'
    'E21' = '2.39%'
    'G21' = '9'
    'D22' = '0.1336'
    'E22' = '2.07%'
    'G22' = '9'
    'D23' = '4.252'
    'E23' = '4.62%'
    'G23' = '9'
    'D24' = '0.1626'
    'E24' = '2.83%'
    'G24' = '9'
    'D25' = '0.001220'
    'E25' = '2.44%'
    'G25' = '9'
    'D26' = '0.004508'
    'E26' = '9.49%'
    'G26' = '9'
    'D27' = '0.0001302'
    'E27' = '8.66%'
    'G27' = '9'
    'D28' = '0.0001775'
    'E28' = '9.79%'
    'G28' = '9'
    'G29' = '9'
    'G30' = '9'
    'G31' = '9'
    'G32' = '9'
    'G33' = '9'
    'G34' = '9'
    'G35' = '9'
    'G36' = '9'
    'G37' = '9'
    'G38' = '9'
    'G39' = '9'
    'D40' = '0.04644'
    'E40' = '8.91%'
    'G40' = '9'
    'D41' = '0.007208'
    'E41' = '7.34%'
    'G41' = '9'
    'D42' = '0.1374'
    'E42' = '10.87%'
    'G42' = '9'
    'D43' = '0.002264'
    'E43' = '2.11%'
    'G43' = '9'
    'D44' = '0.01441'
    'E44' = '19.91%'
    'G44' = '9'
    'D45' = '0.00006200'
    'E45' = '8.90%'
    'G45' = '9'
    'G46' = '9'
    'D47' = '0.01299'
    'E47' = '-0.49%'
    'G47' = '9'
    'G48' = '9'
    'G49' = '9'
    'G50' = '9'
    'G51' = '9'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
